# feat: add 2022-Q1 data
#
# The former "总计" (Total) sheet becomes the new "2022-Q1" quarterly sheet
# (same sheetId/position), and a brand-new "总计" sheet is appended after it
# that repeats the previous total rows plus a new 2022-Q1 row on top.

$wb = $excel.ActiveWorkbook

# Style donors: an existing, untouched quarter sheet already carries the
# workbook's "header" style (bold/centered/bordered) on its header row and
# its bold/bordered look on the leading index column - reuse those via
# copy/paste-special so we mint no new styles.
$styleDonor = $wb.Worksheets.Item(4)
$headerStyleSrc = $styleDonor.Range("B1")
$indexStyleSrc = $styleDonor.Range("A2")

# ---------------------------------------------------------------------
# Step 1: the existing "总计" sheet (5th tab) turns into "2022-Q1" and is
# repopulated with per-fund holdings data, just like the other quarter tabs.
# ---------------------------------------------------------------------
$wsQ1 = $wb.Worksheets.Item(5)
$wsQ1.Cells.Clear()
$wsQ1.Name = "2022-Q1"

$headersQ = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($i = 0; $i -lt $headersQ.Length; $i++) {
    $c = $wsQ1.Cells.Item(1, $i + 2)
    $c.Value = $headersQ[$i]
    $headerStyleSrc.Copy()
    $c.PasteSpecial(-4122)
}

$fundRows = @(
    @("000369", "广发全球医疗保健(QDII) - 人民币", "2.46", "81.85", "1.72", "0.0423", 10),
    @("000370", "广发全球医疗保健(QDII) - 美元", "2.46", "81.85", "1.72", "0.0423", 10),
    @("001092", "广发纳斯达克生物科技指数(QDII)（人民币）", "1.34", "82.00", "3.09", "0.0414", 7),
    @("001093", "广发纳斯达克生物科技指数(QDII)（美元）", "1.34", "82.00", "3.09", "0.0414", 7),
    @("010343", "华宝英国富时100指数（QDII）A", "0.20", "93.65", "7.49", "0.0150", 2),
    @("010344", "华宝英国富时100指数（QDII）C", "0.06", "93.65", "7.49", "0.0045", 2)
)

for ($r = 0; $r -lt $fundRows.Length; $r++) {
    $row = $fundRows[$r]
    $excelRow = $r + 2

    $aCell = $wsQ1.Cells.Item($excelRow, 1)
    $aCell.Value = $r
    $indexStyleSrc.Copy()
    $aCell.PasteSpecial(-4122)

    $bCell = $wsQ1.Cells.Item($excelRow, 2)
    $bCell.NumberFormat = "@"
    $bCell.Value = $row[0]

    $wsQ1.Cells.Item($excelRow, 3).Value = $row[1]

    $dCell = $wsQ1.Cells.Item($excelRow, 4)
    $dCell.NumberFormat = "@"
    $dCell.Value = $row[2]

    $eCell = $wsQ1.Cells.Item($excelRow, 5)
    $eCell.NumberFormat = "@"
    $eCell.Value = $row[3]

    $fCell = $wsQ1.Cells.Item($excelRow, 6)
    $fCell.NumberFormat = "@"
    $fCell.Value = $row[4]

    $gCell = $wsQ1.Cells.Item($excelRow, 7)
    $gCell.NumberFormat = "@"
    $gCell.Value = $row[5]

    $wsQ1.Cells.Item($excelRow, 8).Value = $row[6]
}

# ---------------------------------------------------------------------
# Step 2: add a fresh "总计" sheet right after "2022-Q1", rebuilding the
# summary table with the new quarter prepended to the previous totals.
# ---------------------------------------------------------------------
$wsTotal = $wb.Worksheets.Add($null, $wsQ1)
$wsTotal.Name = "总计"

$headersT = @("日期", "持有数量(只)", "持有市值(亿元)")
for ($i = 0; $i -lt $headersT.Length; $i++) {
    $c = $wsTotal.Cells.Item(1, $i + 2)
    $c.Value = $headersT[$i]
    $headerStyleSrc.Copy()
    $c.PasteSpecial(-4122)
}

$totalRows = @(
    @("2022-Q1", 6, 0.19),
    @("2021-Q4", 2, 0.02),
    @("2021-Q3", 4, 0.11),
    @("2021-Q2", 4, 0.06),
    @("2021-Q1", 4, 0.05)
)

for ($r = 0; $r -lt $totalRows.Length; $r++) {
    $row = $totalRows[$r]
    $excelRow = $r + 2

    $aCell = $wsTotal.Cells.Item($excelRow, 1)
    $aCell.Value = $r
    $indexStyleSrc.Copy()
    $aCell.PasteSpecial(-4122)

    $wsTotal.Cells.Item($excelRow, 2).Value = $row[0]
    $wsTotal.Cells.Item($excelRow, 3).Value = $row[1]
    $wsTotal.Cells.Item($excelRow, 4).Value = $row[2]
}
